$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "2025-04-28 21:29:10"
$ws.Range("B20").Value = 38
